$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.192995548248291
$ws.Range("E2").Value = 731.8656850743482
$ws.Range("F2").Value = 0.02498864556530189
$ws.Range("G2").Value = 0.02194695219327607
$ws.Range("H2").Value = 0.01943285021106969
$ws.Range("I2").Value = 0.01923470042991602
$ws.Range("J2").Value = 0.01794443172541533
$ws.Range("K2").Value = 0.01734847014658198
$ws.Range("L2").Value = 0.01635935181110744
$ws.Range("M2").Value = 0.01576643587777543
$ws.Range("N2").Value = 0.01559840031589149
$ws.Range("O2").Value = 0.01516226709881927
$ws.Range("P2").Value = 0.0150553378218259
$ws.Range("Q2").Value = 0.01479826179918501
$ws.Range("R2").Value = 0.01457490868164615
$ws.Range("S2").Value = 0.01457490868164615
$ws.Range("T2").Value = 0.01456863998132421
$ws.Range("U2").Value = 0.01447322903898607
$ws.Range("V2").Value = 0.01445364772670935
$ws.Range("W2").Value = 0.01433110526677736
$ws.Range("X2").Value = 0.01429847660620753
$ws.Range("Y2").Value = 0.01426638762328164
$ws.Range("C3").Value = 1.224028587341309
$ws.Range("E3").Value = 730.3158482609979
$ws.Range("F3").Value = 0.02485069369717733
$ws.Range("G3").Value = 0.02138248973729333
$ws.Range("H3").Value = 0.01967757464729226
$ws.Range("I3").Value = 0.01898961368953497
$ws.Range("J3").Value = 0.01811619928885757
$ws.Range("K3").Value = 0.01783027823954008
$ws.Range("L3").Value = 0.01675131286624188
$ws.Range("M3").Value = 0.01620644792344667
$ws.Range("N3").Value = 0.0159949578815297
$ws.Range("O3").Value = 0.01570263259504949
$ws.Range("P3").Value = 0.01527896453524441
$ws.Range("Q3").Value = 0.01527896453524441
$ws.Range("R3").Value = 0.01507221530459396
$ws.Range("S3").Value = 0.01485414014138568
$ws.Range("T3").Value = 0.01469557915462234
$ws.Range("U3").Value = 0.01468558598514464
$ws.Range("V3").Value = 0.01448584208913443
$ws.Range("W3").Value = 0.01431956872585758
$ws.Range("X3").Value = 0.01431332623905991
$ws.Range("Y3").Value = 0.01423617637935668
$ws.Range("C4").Value = 1.356996536254883
$ws.Range("E4").Value = 704.4705420291648
$ws.Range("F4").Value = 0.0250579384199764
$ws.Range("G4").Value = 0.02195690500111553
$ws.Range("H4").Value = 0.01959495072566492
$ws.Range("I4").Value = 0.0184839592767913
$ws.Range("J4").Value = 0.0175903625552576
$ws.Range("K4").Value = 0.01711701485722366
$ws.Range("L4").Value = 0.01650247057267285
$ws.Range("M4").Value = 0.01575889368487766
$ws.Range("N4").Value = 0.01551635296109151
$ws.Range("O4").Value = 0.0151210483276648
$ws.Range("P4").Value = 0.01466235754574055
$ws.Range("Q4").Value = 0.01462949826529896
$ws.Range("R4").Value = 0.01449683322082924
$ws.Range("S4").Value = 0.01424868847485683
$ws.Range("T4").Value = 0.01407171880458854
$ws.Range("U4").Value = 0.01407171880458854
$ws.Range("V4").Value = 0.01397295526906087
$ws.Range("W4").Value = 0.01390937022541233
$ws.Range("X4").Value = 0.01377854931563964
$ws.Range("Y4").Value = 0.01373236924033459
$ws.Range("C5").Value = 1.213000059127808
$ws.Range("E5").Value = 732.52077144334
$ws.Range("F5").Value = 0.02517636127337527
$ws.Range("G5").Value = 0.02219450631769508
$ws.Range("H5").Value = 0.02088022884012612
$ws.Range("I5").Value = 0.01930571594687849
$ws.Range("J5").Value = 0.01820560030919771
$ws.Range("K5").Value = 0.01732855869333189
$ws.Range("L5").Value = 0.01732855869333189
$ws.Range("M5").Value = 0.01580124329072133
$ws.Range("N5").Value = 0.01579559626417818
$ws.Range("O5").Value = 0.01568976161701368
$ws.Range("P5").Value = 0.01561773037275083
$ws.Range("Q5").Value = 0.01531888298685674
$ws.Range("R5").Value = 0.01530003496743054
$ws.Range("S5").Value = 0.01494854500320816
$ws.Range("T5").Value = 0.01479752907869003
$ws.Range("U5").Value = 0.01461499495378399
$ws.Range("V5").Value = 0.01455304832516502
$ws.Range("W5").Value = 0.01440176768161375
$ws.Range("X5").Value = 0.01429800914701965
$ws.Range("Y5").Value = 0.0142791573380768
$ws.Range("C6").Value = 1.18199610710144
$ws.Range("E6").Value = 718.5821435960315
$ws.Range("F6").Value = 0.02487803379383803
$ws.Range("G6").Value = 0.02202643385855739
$ws.Range("H6").Value = 0.02005557514555638
$ws.Range("I6").Value = 0.01917533730051952
$ws.Range("J6").Value = 0.01818308221352864
$ws.Range("K6").Value = 0.01733813731047677
$ws.Range("L6").Value = 0.01664827598744452
$ws.Range("M6").Value = 0.01590022561246267
$ws.Range("N6").Value = 0.01585477850565765
$ws.Range("O6").Value = 0.0154719574856087
$ws.Range("P6").Value = 0.01508859191933925
$ws.Range("Q6").Value = 0.01507029995913987
$ws.Range("R6").Value = 0.01490388392549773
$ws.Range("S6").Value = 0.01470896233552659
$ws.Range("T6").Value = 0.014572573309832
$ws.Range("U6").Value = 0.01432252345289496
$ws.Range("V6").Value = 0.01423230642111649
$ws.Range("W6").Value = 0.01414788692306109
$ws.Range("X6").Value = 0.01405105593457219
$ws.Range("Y6").Value = 0.01400744919290509
$ws.Range("C7").Value = 1.243043899536133
$ws.Range("E7").Value = 720.6989829430386
$ws.Range("F7").Value = 0.02500324008744786
$ws.Range("G7").Value = 0.02171208522362391
$ws.Range("H7").Value = 0.01891917729849667
$ws.Range("I7").Value = 0.01857565267749674
$ws.Range("J7").Value = 0.01717513031805876
$ws.Range("K7").Value = 0.01717513031805876
$ws.Range("L7").Value = 0.01622743093615392
$ws.Range("M7").Value = 0.0156561498309973
$ws.Range("N7").Value = 0.01492983308180764
$ws.Range("O7").Value = 0.01492983308180764
$ws.Range("P7").Value = 0.01492983308180764
$ws.Range("Q7").Value = 0.01478320019913647
$ws.Range("R7").Value = 0.01447594598501875
$ws.Range("S7").Value = 0.01438294076666856
$ws.Range("T7").Value = 0.01437318440031735
$ws.Range("U7").Value = 0.01437219420832749
$ws.Range("V7").Value = 0.01437219420832749
$ws.Range("W7").Value = 0.01423772058458032
$ws.Range("X7").Value = 0.0141253217996441
$ws.Range("Y7").Value = 0.01404871311779802
$ws.Range("C8").Value = 1.125962257385254
$ws.Range("E8").Value = 725.9503162243218
$ws.Range("F8").Value = 0.02565526101116538
$ws.Range("G8").Value = 0.02142918479130106
$ws.Range("H8").Value = 0.01987471020284703
$ws.Range("I8").Value = 0.01885862481892634
$ws.Range("J8").Value = 0.01749663321410973
$ws.Range("K8").Value = 0.01721769175635604
$ws.Range("L8").Value = 0.01665135011585184
$ws.Range("M8").Value = 0.01621726467583233
$ws.Range("N8").Value = 0.01621726467583233
$ws.Range("O8").Value = 0.0154452800018192
$ws.Range("P8").Value = 0.01533788327949087
$ws.Range("Q8").Value = 0.01507702493664742
$ws.Range("R8").Value = 0.0148995861537669
$ws.Range("S8").Value = 0.01466025764726464
$ws.Range("T8").Value = 0.01466025764726464
$ws.Range("U8").Value = 0.01455537825373848
$ws.Range("V8").Value = 0.01432183647550558
$ws.Range("W8").Value = 0.01424394963428683
$ws.Range("X8").Value = 0.01418002715359705
$ws.Range("Y8").Value = 0.01415107828897313
$ws.Range("C9").Value = 1.104996919631958
$ws.Range("E9").Value = 753.2128301522007
$ws.Range("F9").Value = 0.02544959897322169
$ws.Range("G9").Value = 0.02217108256278091
$ws.Range("H9").Value = 0.02014027365615055
$ws.Range("I9").Value = 0.0190137965797641
$ws.Range("J9").Value = 0.01890171115846173
$ws.Range("K9").Value = 0.01742171362875708
$ws.Range("L9").Value = 0.017080220748544
$ws.Range("M9").Value = 0.01608815032188172
$ws.Range("N9").Value = 0.01608815032188172
$ws.Range("O9").Value = 0.01582658847494658
$ws.Range("P9").Value = 0.01579846187927887
$ws.Range("Q9").Value = 0.01548992952066672
$ws.Range("R9").Value = 0.01548992952066672
$ws.Range("S9").Value = 0.0153517140403353
$ws.Range("T9").Value = 0.01517179200735323
$ws.Range("U9").Value = 0.01485674339299354
$ws.Range("V9").Value = 0.01485674339299354
$ws.Range("W9").Value = 0.01480150057749943
$ws.Range("X9").Value = 0.01472443309641193
$ws.Range("Y9").Value = 0.01468251130900976
$ws.Range("C10").Value = 1.274001598358154
$ws.Range("E10").Value = 735.5095723066861
$ws.Range("F10").Value = 0.0248661752498197
$ws.Range("G10").Value = 0.02104954793624803
$ws.Range("H10").Value = 0.01933114715552968
$ws.Range("I10").Value = 0.01799173343314044
$ws.Range("J10").Value = 0.01784957692807966
$ws.Range("K10").Value = 0.01690344968458247
$ws.Range("L10").Value = 0.01651886557780286
$ws.Range("M10").Value = 0.01585943608549592
$ws.Range("N10").Value = 0.01585943608549592
$ws.Range("O10").Value = 0.0156407934300766
$ws.Range("P10").Value = 0.01555555273406623
$ws.Range("Q10").Value = 0.01505762722664242
$ws.Range("R10").Value = 0.01488702045274615
$ws.Range("S10").Value = 0.01488702045274615
$ws.Range("T10").Value = 0.01474253436415187
$ws.Range("U10").Value = 0.01471447362005528
$ws.Range("V10").Value = 0.01460734539775811
$ws.Range("W10").Value = 0.01439545235856393
$ws.Range("X10").Value = 0.01439545235856393
$ws.Range("Y10").Value = 0.01433741856348316
$ws.Range("C11").Value = 1.044998407363892
$ws.Range("E11").Value = 724.639497773429
$ws.Range("F11").Value = 0.02563122447068902
$ws.Range("G11").Value = 0.02221389097834603
$ws.Range("H11").Value = 0.01977291570025325
$ws.Range("I11").Value = 0.01860012697132081
$ws.Range("J11").Value = 0.01818854816680721
$ws.Range("K11").Value = 0.01730812771839519
$ws.Range("L11").Value = 0.01691892478077543
$ws.Range("M11").Value = 0.01641038178983463
$ws.Range("N11").Value = 0.01571405475518844
$ws.Range("O11").Value = 0.01568725221936229
$ws.Range("P11").Value = 0.01542743328736361
$ws.Range("Q11").Value = 0.01495031502956067
$ws.Range("R11").Value = 0.0148251271002158
$ws.Range("S11").Value = 0.0146737571664184
$ws.Range("T11").Value = 0.014567572552859
$ws.Range("U11").Value = 0.01436357730645542
$ws.Range("V11").Value = 0.01430008063248147
$ws.Range("W11").Value = 0.0142477258429142
$ws.Range("X11").Value = 0.01417948397815851
$ws.Range("Y11").Value = 0.01412552627238653
